$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new drawback entries in column B (rows 3 and 4)
# Set B4 first so the shared-string table indexes "Attacks are simple to
# carry out" before "Large number of attack vectors", matching the target.
$ws.Range("B4").Value = "Attacks are simple to carry out"
$ws.Range("B3").Value = "Large number of attack vectors"

# Update the active selection to match the new cursor position
$ws.Range("B7").Select()
